# Auto-generated update of market-price derived columns (H:N) across all item sheets.
# Mirrors a scheduled market-data refresh: currentAveragePrice / LevePrice / LeveProfit
# columns are recomputed values with no underlying formulas in this workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 1140.0625
$ws.Range("J17").Value = 928.93616
$ws.Range("L17").Value = 2786.80848
$ws.Range("N17").Value = -3122.80848
$ws.Range("H40").Value = 1759.1538
$ws.Range("I40").Value = 2056
$ws.Range("J40").Value = 1573.625
$ws.Range("K40").Value = 2056
$ws.Range("L40").Value = 1573.625
$ws.Range("M40").Value = -1881
$ws.Range("N40").Value = -1923.625
$ws.Range("H55").Value = 291.92856
$ws.Range("J55").Value = 354.14285
$ws.Range("L55").Value = 354.14285
$ws.Range("N55").Value = -782.14285
$ws.Range("H99").Value = 953.8
$ws.Range("I99").Value = 327
$ws.Range("J99").Value = 1894
$ws.Range("K99").Value = 981
$ws.Range("L99").Value = 5682
$ws.Range("M99").Value = 517
$ws.Range("N99").Value = -8678
$ws.Range("H101").Value = 1199.5
$ws.Range("J101").Value = 1999
$ws.Range("L101").Value = 5997
$ws.Range("N101").Value = -9241
$ws.Range("H111").Value = 2966.6667
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2966.6667
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 8900.000100000001
$ws.Range("N111").Value = -15034.0001
$ws.Range("M111").ClearContents()
$ws.Range("H132").Value = 1398.9286
$ws.Range("I132").Value = 1391.1538
$ws.Range("K132").Value = 4173.4614
$ws.Range("M132").Value = -1643.4614
$ws.Range("H135").Value = 444.41177
$ws.Range("I135").Value = 444.41177
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3999.70593
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -1464.70593
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2117
$ws.Range("I137").Value = 1675.5
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 5026.5
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -2476.5
$ws.Range("N137").Value = -14100
$ws.Range("H141").Value = 2158320.5
$ws.Range("I141").Value = 4003440.2
$ws.Range("K141").Value = 12010320.6
$ws.Range("M141").Value = -12005140.6

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 3835.7646
$ws.Range("I32").Value = 2343.907
$ws.Range("K32").Value = 2343.907
$ws.Range("M32").Value = -2056.907
$ws.Range("H61").Value = 2474.681
$ws.Range("I61").Value = 1991.4651
$ws.Range("J61").Value = 7669.25
$ws.Range("K61").Value = 1991.4651
$ws.Range("L61").Value = 7669.25
$ws.Range("M61").Value = -1779.4651
$ws.Range("N61").Value = -8093.25
$ws.Range("H88").Value = 21790.545
$ws.Range("J88").Value = 28974.625
$ws.Range("L88").Value = 28974.625
$ws.Range("N88").Value = -29786.625
$ws.Range("H91").Value = 21790.545
$ws.Range("J91").Value = 28974.625
$ws.Range("L91").Value = 28974.625
$ws.Range("N91").Value = -31782.625
$ws.Range("H136").Value = 2474.681
$ws.Range("I136").Value = 1991.4651
$ws.Range("J136").Value = 7669.25
$ws.Range("K136").Value = 5974.3953
$ws.Range("L136").Value = 23007.75
$ws.Range("M136").Value = -3424.3953
$ws.Range("N136").Value = -28107.75

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -277
$ws.Range("N22").Value = -946
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 97402.83
$ws.Range("I86").Value = 1734.75
$ws.Range("K86").Value = 1734.75
$ws.Range("M86").Value = -611.75
$ws.Range("H89").Value = 97402.83
$ws.Range("I89").Value = 1734.75
$ws.Range("K89").Value = 8673.75
$ws.Range("M89").Value = -3057.75
$ws.Range("H94").Value = 590.28
$ws.Range("I94").Value = 557.4545000000001
$ws.Range("K94").Value = 557.4545000000001
$ws.Range("M94").Value = -106.4545000000001

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H19").Value = 1084.3334
$ws.Range("I19").Value = 751.5
$ws.Range("J19").Value = 1750
$ws.Range("K19").Value = 751.5
$ws.Range("L19").Value = 1750
$ws.Range("M19").Value = -581.5
$ws.Range("N19").Value = -2090
$ws.Range("H22").Value = 1220
$ws.Range("H24").Value = 1084.3334
$ws.Range("I24").Value = 751.5
$ws.Range("J24").Value = 1750
$ws.Range("K24").Value = 751.5
$ws.Range("L24").Value = 1750
$ws.Range("M24").Value = -581.5
$ws.Range("N24").Value = -2090
$ws.Range("H31").Value = 2380.4897
$ws.Range("I31").Value = 1574.3823
$ws.Range("J31").Value = 4207.6665
$ws.Range("K31").Value = 1574.3823
$ws.Range("L31").Value = 4207.6665
$ws.Range("M31").Value = -1279.3823
$ws.Range("N31").Value = -4797.6665
$ws.Range("H34").Value = 2380.4897
$ws.Range("I34").Value = 1574.3823
$ws.Range("J34").Value = 4207.6665
$ws.Range("K34").Value = 1574.3823
$ws.Range("L34").Value = 4207.6665
$ws.Range("M34").Value = -1372.3823
$ws.Range("N34").Value = -4611.6665
$ws.Range("H58").Value = 1554319.8
$ws.Range("I58").Value = 3345559.8
$ws.Range("J58").Value = 1911.7333
$ws.Range("K58").Value = 3345559.8
$ws.Range("L58").Value = 1911.7333
$ws.Range("M58").Value = -3345356.8
$ws.Range("N58").Value = -2317.7333
$ws.Range("H132").Value = 1889.7715
$ws.Range("I132").Value = 1138.6818
$ws.Range("K132").Value = 3416.0454
$ws.Range("M132").Value = -886.0454
$ws.Range("H136").Value = 1554319.8
$ws.Range("I136").Value = 3345559.8
$ws.Range("J136").Value = 1911.7333
$ws.Range("K136").Value = 10036679.4
$ws.Range("L136").Value = 5735.199900000001
$ws.Range("M136").Value = -10034129.4
$ws.Range("N136").Value = -10835.1999

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H56").Value = 7960.25
$ws.Range("I56").Value = 7960.25
$ws.Range("K56").Value = 7960.25
$ws.Range("M56").Value = -7430.25
$ws.Range("H107").Value = 906.3333
$ws.Range("J107").Value = 906.3333
$ws.Range("L107").Value = 2718.9999
$ws.Range("N107").Value = -6558.9999
$ws.Range("H131").Value = 9626.089
$ws.Range("J131").Value = 10503.292
$ws.Range("L131").Value = 31509.876
$ws.Range("N131").Value = -41589.876

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 2633
$ws.Range("I80").Value = 2649.5
$ws.Range("K80").Value = 2649.5
$ws.Range("M80").Value = -1651.5
$ws.Range("H83").Value = 2633
$ws.Range("I83").Value = 2649.5
$ws.Range("K83").Value = 13247.5
$ws.Range("M83").Value = -8255.5
$ws.Range("H102").Value = 3291.6924
$ws.Range("I102").Value = 4969.75
$ws.Range("J102").Value = 2545.889
$ws.Range("K102").Value = 4969.75
$ws.Range("L102").Value = 2545.889
$ws.Range("M102").Value = -3347.75
$ws.Range("N102").Value = -5789.889

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 2521.7778
$ws.Range("J22").Value = 1713.7142
$ws.Range("L22").Value = 1713.7142
$ws.Range("N22").Value = -2303.7142
$ws.Range("H27").Value = 2521.7778
$ws.Range("J27").Value = 1713.7142
$ws.Range("L27").Value = 1713.7142
$ws.Range("N27").Value = -1927.7142
$ws.Range("H46").Value = 2078.2856
$ws.Range("I46").Value = 1099.8
$ws.Range("K46").Value = 1099.8
$ws.Range("M46").Value = -911.8
$ws.Range("H61").Value = 5874.75
$ws.Range("I61").Value = 6166.3335
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 6166.3335
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -5964.3335
$ws.Range("N61").Value = -5404
$ws.Range("H113").Value = 5874.75
$ws.Range("I113").Value = 6166.3335
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 6166.3335
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -3996.3335
$ws.Range("N113").Value = -9340
$ws.Range("H122").Value = 6179.1113
$ws.Range("I122").Value = 1803.5
$ws.Range("J122").Value = 7429.2856
$ws.Range("K122").Value = 5410.5
$ws.Range("L122").Value = 22287.8568
$ws.Range("M122").Value = -2960.5
$ws.Range("N122").Value = -27187.8568
$ws.Range("H132").Value = 1699.0555
$ws.Range("I132").Value = 1479.4
$ws.Range("J132").Value = 1973.625
$ws.Range("K132").Value = 4438.200000000001
$ws.Range("L132").Value = 5920.875
$ws.Range("M132").Value = -1908.200000000001
$ws.Range("N132").Value = -10980.875

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H126").Value = 8592.5625
$ws.Range("J126").Value = 9139.714
$ws.Range("L126").Value = 27419.142
$ws.Range("N126").Value = -32359.142
$ws.Range("H132").Value = 1894.963
$ws.Range("I132").Value = 1508.7368
$ws.Range("J132").Value = 2812.25
$ws.Range("K132").Value = 4526.2104
$ws.Range("L132").Value = 8436.75
$ws.Range("M132").Value = -1996.2104
$ws.Range("N132").Value = -13496.75
$ws.Range("H136").Value = 15016988
$ws.Range("I136").Value = 20577838
$ws.Range("K136").Value = 61733514
$ws.Range("M136").Value = -61730964
